$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert "Jalap" (id 367, severity 85) right before the existing "Jalla" row (alphabetical order).
$ws.Rows.Item(26).Insert()
$ws.Range("A26").Value = 367
$ws.Range("B26").Value = "Jalap"
$ws.Range("C26").Value = 85

# Insert "kashanda" (id 366, severity 25) right before the existing "ko't" row.
$ws.Rows.Item(130).Insert()
$ws.Range("A130").Value = 366
$ws.Range("B130").Value = "kashanda"
$ws.Range("C130").Value = 25

# Insert "oneni ami" (id 365, severity 100) right before the existing "oom" row.
$ws.Rows.Item(142).Insert()
$ws.Range("A142").Value = 365
$ws.Range("B142").Value = "oneni ami"
$ws.Range("C142").Value = 100
